$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 847, shifting existing rows 847:950 down to 848:951
$ws.Rows.Item(847).Insert()

# Populate the new row 847 with the new record.
# Columns A,B,C,E,F,G,H,N,Q,R carry the same constant values as every other
# row in this data block, so copy them down from the row above (846).
$ws.Range("A847").Value = $ws.Range("A846").Value2
$ws.Range("B847").Value = $ws.Range("B846").Value2
$ws.Range("C847").Value = $ws.Range("C846").Value2
$ws.Range("D847").Value = 45124
$ws.Range("E847").Value = $ws.Range("E846").Value2
$ws.Range("F847").Value = $ws.Range("F846").Value2
$ws.Range("G847").Value = $ws.Range("G846").Value2
$ws.Range("H847").Value = $ws.Range("H846").Value2
$ws.Range("I847").Value = "1a (guarda)"
$ws.Range("J847").Value = 1800
$ws.Range("K847").Value = 600
$ws.Range("L847").Value = 700
$ws.Range("M847").Value = 650
$ws.Range("N847").Value = $ws.Range("N846").Value2
$ws.Range("O847").Value = "Región de O'Higgins"
$ws.Range("P847").Value = 650
$ws.Range("Q847").Value = $ws.Range("Q846").Value2
$ws.Range("R847").Value = $ws.Range("R846").Value2
